$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert "Min Voltage" column before the existing
#     Price/Quantity/Total columns (F,G,H -> G,H,I), shifting only
#     those three header cells (N1/O1 and R4:R9 must stay put).
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = $ws.Range("H1").Value()
$ws.Range("H1").Value = $ws.Range("G1").Value()
$ws.Range("G1").Value = $ws.Range("F1").Value()
$ws.Range("F1").Value = "Min Voltage"

# New column F width
$ws.Columns.Item(6).ColumnWidth = 12.3

# --- Rows 6-9 (the Electrolytic-cap rows that are no longer relevant
#     now capacitor selection is finished) become hidden.
$ws.Rows.Item(6).Hidden = $true
$ws.Rows.Item(7).Hidden = $true
$ws.Rows.Item(8).Hidden = $true
$ws.Rows.Item(9).Hidden = $true

# --- Fill in the now-selected ceramic capacitors (rows 10-12).
$ws.Range("D11").Value = "GRM21BR61H106KE43L (at JLC)"
$ws.Range("C11").Value = "10uF Ceramic Cap"
$ws.Range("C10").Value = "22nf Ceramic Cap"
$ws.Range("E11").Value = "0805(imperial)"
$ws.Range("B11").Value = "C8,C9,C12,C13,C22,C54,C57"
$ws.Range("E10").Value = "0603(imperial)"
$ws.Range("C12").Value = "2.2uF Ceramic Cap"
$ws.Range("D10").Value = "CL10B223KB8NNNC (JLC)"
$ws.Range("D12").Value = "0805F225M500NT (JLC)"
$ws.Range("E12").Value = "0805(imperial)"
$ws.Range("B12").Value = "C16"
$ws.Range("B10").Value = "CTTC1"

# Match the author's final selection.
$ws.Range("E12").Select()
